# Apply updated cryptocurrency price/volume data to the worksheet,
# matching the source CSV refresh performed by the scheduled GitHub Action.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '68.631.43'
$cell.ClearFormats()
$ws.Range("E2").Value = '  -0.66%  '

# Row 3
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '3.746.63'
$cell.ClearFormats()
$ws.Range("E3").Value = '  -1.87%  '

# Row 4
$ws.Range("E4").Value = '  +0.15%  '

# Row 5
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '625.44'
$cell.ClearFormats()
$ws.Range("E5").Value = '  -0.70%  '

# Row 6
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '164.57'
$cell.ClearFormats()
$ws.Range("E6").Value = '  -0.46%  '

# Row 7
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '3.745.17'
$cell.ClearFormats()
$ws.Range("E7").Value = '  -1.85%  '

# Row 8
$ws.Range("E8").Value = '  -0.02%  '

# Row 9
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.520'
$cell.ClearFormats()
$ws.Range("E9").Value = '  +0.07%  '

# Row 10
$ws.Range("E10").Value = '  -2.87%  '

# Row 11
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.453'
$cell.ClearFormats()
$ws.Range("E11").Value = '  -0.17%  '

# Row 12
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '6.88'
$cell.ClearFormats()
$ws.Range("E12").Value = '  +3.80%  '

# Row 13
$ws.Range("E13").Value = '  -5.42%  '

# Row 14
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '34.59'
$cell.ClearFormats()
$ws.Range("E14").Value = '  -3.95%  '

# Row 15
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '4.379.22'
$cell.ClearFormats()
$ws.Range("E15").Value = '  -1.76%  '

# Row 16
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '3.741.45'
$cell.ClearFormats()
$ws.Range("E16").Value = '  -1.22%  '

# Row 17
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '68.619.33'
$cell.ClearFormats()
$ws.Range("E17").Value = '  -0.63%  '

# Row 18
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '17.52'
$cell.ClearFormats()
$ws.Range("E18").Value = '  -2.87%  '

# Row 19
$ws.Range("E19").Value = '  -0.48%  '

# Row 20
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '6.95'
$cell.ClearFormats()
$ws.Range("E20").Value = '  -2.33%  '

# Row 21
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '467.11'
$cell.ClearFormats()
$ws.Range("E21").Value = '  +0.30%  '

# Row 22
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '9.42'
$cell.ClearFormats()
$ws.Range("E22").Value = '  -2.32%  '

# Row 23
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '0.698'
$cell.ClearFormats()
$ws.Range("E23").Value = '  -1.37%  '

# Row 24
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '81.27'
$cell.ClearFormats()
$ws.Range("E24").Value = '  -2.95%  '

# Row 25
$ws.Range("E25").Value = '  -7.05%  '

# Row 26
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '12.04'
$cell.ClearFormats()
$ws.Range("E26").Value = '  +0.67%  '

# Row 27
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '2.08'
$cell.ClearFormats()
$ws.Range("E27").Value = '  -3.18%  '

# Row 28
$ws.Range("B28").Value = 'Dai'
$ws.Range("C28").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.ClearFormats()
$ws.Range("E28").Value = '  -0.06%  '

# Row 29
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '10.00'
$cell.ClearFormats()
$ws.Range("E29").Value = '  -0.33%  '

# Row 30
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '3.894.88'
$cell.ClearFormats()
$ws.Range("E30").Value = '  -1.82%  '

# Row 31
$ws.Range("E31").Value = '  +1.14%  '

# Row 32
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '2.64'
$cell.ClearFormats()
$ws.Range("E32").Value = '  -1.79%  '

# Row 33
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '7.07'
$cell.ClearFormats()
$ws.Range("E33").Value = '  -2.90%  '

# Row 34
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '0.176'
$cell.ClearFormats()
$ws.Range("E34").Value = '  +17.89%  '

# Row 35
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '28.22'
$cell.ClearFormats()
$ws.Range("E35").Value = '  -3.31%  '

# Row 36
$ws.Range("E36").Value = '  -0.11%  '

# Row 37
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '3.701.13'
$cell.ClearFormats()
$ws.Range("E37").Value = '  -1.62%  '

# Row 38
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '8.82'
$cell.ClearFormats()
$ws.Range("E38").Value = '  -2.94%  '

# Row 39
$ws.Range("E39").Value = '  -1.43%  '

# Row 40
$ws.Range("B40").Value = 'Filecoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '5.72'
$cell.ClearFormats()
$ws.Range("E40").Value = '  -3.15%  '

# Row 41
$ws.Range("B41").Value = 'dogwifhat'
$ws.Range("C41").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '3.21'
$cell.ClearFormats()
$ws.Range("E41").Value = '  -6.06%  '

# Row 42
$ws.Range("E42").Value = '  +0.01%  '

# Row 43
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '0.951'
$cell.ClearFormats()
$ws.Range("E43").Value = '  -2.99%  '

# Row 44
$ws.Range("E44").Value = '  -0.01%  '

# Row 45
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '44.12'
$cell.ClearFormats()
$ws.Range("E45").Value = '  +3.93%  '

# Row 46
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '156.34'
$cell.ClearFormats()
$ws.Range("E46").Value = '  -0.83%  '

# Row 47
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '1.94'
$cell.ClearFormats()
$ws.Range("E47").Value = '  +2.23%  '

# Row 48
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '47.23'
$cell.ClearFormats()
$ws.Range("E48").Value = '  +0.74%  '

# Row 49
$ws.Range("E49").Value = '  -3.50%  '

# Row 50
$ws.Range("E50").Value = '  -2.80%  '

# Row 51
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '8.31'
$cell.ClearFormats()
$ws.Range("E51").Value = '  -1.67%  '
